$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 92
$ws.Range("F2").Value = 64
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 70
$ws.Range("E3").Value = 37
$ws.Range("F3").Value = 28
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 29
$ws.Range("E4").Value = 43
$ws.Range("F4").Value = 27
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = 39
$ws.Range("E5").Value = 133
$ws.Range("F5").Value = 90
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 101
$ws.Range("E6").Value = 43
$ws.Range("F6").Value = 31
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 41
$ws.Range("E7").Value = 30
$ws.Range("F7").Value = 16
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 21
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 6
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 5
$ws.Range("E10").Value = 510
$ws.Range("F10").Value = 256
$ws.Range("G10").Value = 95
$ws.Range("H10").Value = 351
$ws.Range("E11").Value = 333
$ws.Range("F11").Value = 186
$ws.Range("G11").Value = 68
$ws.Range("H11").Value = 254
$ws.Range("E12").Value = 499
$ws.Range("F12").Value = 273
$ws.Range("G12").Value = 83
$ws.Range("H12").Value = 356
$ws.Range("E13").Value = 126
$ws.Range("F13").Value = 68
$ws.Range("G13").Value = 34
$ws.Range("H13").Value = 102
$ws.Range("E14").Value = 124
$ws.Range("F14").Value = 66
$ws.Range("G14").Value = 35
$ws.Range("H14").Value = 101
$ws.Range("E15").Value = 159
$ws.Range("F15").Value = 70
$ws.Range("G15").Value = 49
$ws.Range("H15").Value = 119
$ws.Range("E16").Value = 196
$ws.Range("F16").Value = 100
$ws.Range("G16").Value = 48
$ws.Range("H16").Value = 148
$ws.Range("E17").Value = 98
$ws.Range("F17").Value = 53
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 77
$ws.Range("E18").Value = 51
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 43
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = 8
$ws.Range("E20").Value = 86
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 37
$ws.Range("H20").Value = 67
$ws.Range("E21").Value = 136
$ws.Range("F21").Value = 74
$ws.Range("G21").Value = 31
$ws.Range("H21").Value = 105
$ws.Range("E22").Value = 164
$ws.Range("F22").Value = 88
$ws.Range("G22").Value = 42
$ws.Range("H22").Value = 130
$ws.Range("E23").Value = 197
$ws.Range("F23").Value = 92
$ws.Range("G23").Value = 51
$ws.Range("H23").Value = 143
$ws.Range("E24").Value = 206
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 30
$ws.Range("H24").Value = 141
$ws.Range("E25").Value = 257
$ws.Range("F25").Value = 125
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = 185
$ws.Range("E26").Value = 152
$ws.Range("F26").Value = 94
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 119
$ws.Range("E27").Value = 317
$ws.Range("F27").Value = 158
$ws.Range("G27").Value = 81
$ws.Range("H27").Value = 239
$ws.Range("E28").Value = 191
$ws.Range("F28").Value = 74
$ws.Range("G28").Value = 52
$ws.Range("H28").Value = 126
$ws.Range("E29").Value = 163
$ws.Range("F29").Value = 92
$ws.Range("G29").Value = 41
$ws.Range("H29").Value = 133
$ws.Range("E30").Value = 199
$ws.Range("F30").Value = 119
$ws.Range("G30").Value = 52
$ws.Range("H30").Value = 171
$ws.Range("E31").Value = 71
$ws.Range("F31").Value = 32
$ws.Range("G31").Value = 28
$ws.Range("H31").Value = 60
$ws.Range("E32").Value = 179
$ws.Range("F32").Value = 106
$ws.Range("G32").Value = 38
$ws.Range("H32").Value = 144
$ws.Range("E33").Value = 283
$ws.Range("F33").Value = 144
$ws.Range("G33").Value = 89
$ws.Range("H33").Value = 233
$ws.Range("E34").Value = 211
$ws.Range("F34").Value = 139
$ws.Range("G34").Value = 39
$ws.Range("H34").Value = 178
$ws.Range("E35").Value = 145
$ws.Range("F35").Value = 90
$ws.Range("G35").Value = 27
$ws.Range("H35").Value = 117
$ws.Range("E36").Value = 68
$ws.Range("F36").Value = 39
$ws.Range("G36").Value = 10
$ws.Range("H36").Value = 49
$ws.Range("E37").Value = 152
$ws.Range("F37").Value = 74
$ws.Range("G37").Value = 37
$ws.Range("H37").Value = 111
$ws.Range("E38").Value = 89
$ws.Range("F38").Value = 55
$ws.Range("G38").Value = 16
$ws.Range("H38").Value = 71
$ws.Range("E39").Value = 178
$ws.Range("F39").Value = 88
$ws.Range("G39").Value = 51
$ws.Range("H39").Value = 139
$ws.Range("E40").Value = 255
$ws.Range("F40").Value = 120
$ws.Range("G40").Value = 80
$ws.Range("H40").Value = 200
$ws.Range("E41").Value = 380
$ws.Range("F41").Value = 177
$ws.Range("G41").Value = 92
$ws.Range("H41").Value = 269
$ws.Range("E42").Value = 370
$ws.Range("F42").Value = 204
$ws.Range("G42").Value = 60
$ws.Range("H42").Value = 264
$ws.Range("E43").Value = 113
$ws.Range("F43").Value = 61
$ws.Range("G43").Value = 28
$ws.Range("H43").Value = 89
$ws.Range("E44").Value = 305
$ws.Range("F44").Value = 154
$ws.Range("G44").Value = 68
$ws.Range("H44").Value = 222
$ws.Range("E45").Value = 139
$ws.Range("F45").Value = 70
$ws.Range("G45").Value = 39
$ws.Range("H45").Value = 109
$ws.Range("E46").Value = 306
$ws.Range("F46").Value = 168
$ws.Range("G46").Value = 63
$ws.Range("H46").Value = 231
$ws.Range("E47").Value = 440
$ws.Range("F47").Value = 221
$ws.Range("G47").Value = 92
$ws.Range("H47").Value = 313
$ws.Range("E48").Value = 199
$ws.Range("F48").Value = 87
$ws.Range("G48").Value = 44
$ws.Range("H48").Value = 131
$ws.Range("E49").Value = 280
$ws.Range("F49").Value = 120
$ws.Range("G49").Value = 87
$ws.Range("H49").Value = 207
$ws.Range("E50").Value = 236
$ws.Range("F50").Value = 111
$ws.Range("G50").Value = 73
$ws.Range("H50").Value = 184
$ws.Range("E51").Value = 226
$ws.Range("F51").Value = 98
$ws.Range("G51").Value = 72
$ws.Range("H51").Value = 170
$ws.Range("E52").Value = 25
$ws.Range("F52").Value = 12
$ws.Range("G52").Value = 8
$ws.Range("H52").Value = 20

Write-Output "Updated rows 2-52 E:H"
